$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update "invalidCredentialData" sheet (sheet1) ---
# Row 3 (was peter/peter123/Danish/Invalid...) -> king/king123/Danish/Invalid...
$ws1.Range("A3").Value = "king"
$ws1.Range("B3").Value = "king123"

# Remove the old rows 4 and 5 (peter12 / 223frrr test cases)
$ws1.Rows.Item(4).EntireRow.Delete()
$ws1.Rows.Item(4).EntireRow.Delete()

# Update the view/selection on sheet1
$null = $ws1.Range("A1:C2").Select()

# --- Add the new "CheckHeaderAndVersionData" sheet right after invalidCredentialData ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "CheckHeaderAndVersionData"

# Header row
$ws2.Range("A1").Value = "User Name"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "Language"
$ws2.Range("D1").Value = "Expected About Header"
$ws2.Range("E1").Value = "Expected Version"

# Data row 2 - About header mismatch case
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pass"
$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("D2").Value = "About OpenEMR787"
$ws2.Range("E2").Value = "Version Number: v6.0.0 (2)"

# Data row 3 - matching case
$ws2.Range("A3").Value = "admin"
$ws2.Range("B3").Value = "pass"
$ws2.Range("C3").Value = "English (Indian)"
$ws2.Range("D3").Value = "About OpenEMR"
$ws2.Range("E3").Value = "Version Number: v6.0.0 (2)"

# Column widths (auto-fit sized in the authored workbook)
$ws2.Columns.Item(1).ColumnWidth = 10.7109375
$ws2.Columns.Item(2).ColumnWidth = 9.42578125
$ws2.Columns.Item(3).ColumnWidth = 14.85546875
$ws2.Columns.Item(4).ColumnWidth = 22.28515625
$ws2.Columns.Item(5).ColumnWidth = 16.5703125

# Selection/active cell on the new sheet
$null = $ws2.Range("D4").Select()
